$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

# Delete row 57 (Caso 7098, UGARTE MANUEL 3484) - the rest of the rows shift up by one
$ws.Rows.Item(57).Delete()
